$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('LP1912')
$updates = @(
    @{Row=2; Col=1; Value='Última actualización: 08:01:10'},
    @{Row=3; Col=1; Value='Total filas: 106'},
    @{Row=26; Col=1; Value='06:25:28'},
    @{Row=26; Col=3; Value='23_HERNANDEZ'},
    @{Row=26; Col=4; Value=5},
    @{Row=27; Col=1; Value='05:55:02'},
    @{Row=27; Col=3; Value='86_EST CHICA-ESC AGRARIA'},
    @{Row=27; Col=4; Value=35},
    @{Row=32; Col=3; Value='215C_EL PATO'},
    @{Row=33; Col=3; Value='14_ABASTO'},
    @{Row=37; Col=1; Value='06:54:06'},
    @{Row=37; Col=3; Value='23_HERNANDEZ'},
    @{Row=37; Col=4; Value=11},
    @{Row=38; Col=1; Value='06:25:28'},
    @{Row=38; Col=3; Value='15_ABASTO'},
    @{Row=38; Col=4; Value=40},
    @{Row=39; Col=3; Value='15_ABASTO'},
    @{Row=40; Col=3; Value='225_GOMEZ'},
    @{Row=67; Col=1; Value='07:17:59'},
    @{Row=67; Col=3; Value='11_ETCHEVERRY'},
    @{Row=67; Col=4; Value=46},
    @{Row=68; Col=1; Value='08:01:10'},
    @{Row=68; Col=3; Value='23_HERNANDEZ'},
    @{Row=68; Col=4; Value=2},
    @{Row=69; Col=1; Value='08:01:10'},
    @{Row=69; Col=4; Value=3},
    @{Row=71; Col=1; Value='07:17:59'},
    @{Row=71; Col=3; Value='15_ABASTO'},
    @{Row=71; Col=4; Value=54},
    @{Row=72; Col=1; Value='08:01:10'},
    @{Row=72; Col=3; Value='16_SANTA ANA'},
    @{Row=72; Col=4; Value=10},
    @{Row=73; Col=1; Value='07:17:59'},
    @{Row=73; Col=3; Value='10_OLMOS'},
    @{Row=73; Col=4; Value=55},
    @{Row=74; Col=1; Value='08:01:10'},
    @{Row=74; Col=3; Value='15_ABASTO'},
    @{Row=74; Col=4; Value=11},
    @{Row=75; Col=1; Value='08:01:10'},
    @{Row=75; Col=4; Value=12},
    @{Row=78; Col=3; Value='16_P MOR-SANTA ANA'},
    @{Row=79; Col=3; Value='215B_EL PATO'},
    @{Row=80; Col=1; Value='08:01:10'},
    @{Row=80; Col=3; Value='16_P MOR-SANTA ANA'},
    @{Row=80; Col=4; Value=22},
    @{Row=81; Col=1; Value='08:01:10'},
    @{Row=81; Col=3; Value='215B_EL PATO'},
    @{Row=81; Col=4; Value=22},
    @{Row=82; Col=1; Value='08:01:10'},
    @{Row=82; Col=2; Value='08:24'},
    @{Row=82; Col=3; Value='14_ABASTO'},
    @{Row=82; Col=4; Value=23},
    @{Row=83; Col=1; Value='07:17:59'},
    @{Row=83; Col=2; Value='08:26'},
    @{Row=83; Col=4; Value=69},
    @{Row=84; Col=1; Value='08:01:10'},
    @{Row=84; Col=2; Value='08:27'},
    @{Row=84; Col=3; Value='84_COLONIA URQUIZA-ESC 49'},
    @{Row=84; Col=4; Value=26},
    @{Row=85; Col=2; Value='08:30'},
    @{Row=85; Col=3; Value='23_HERNANDEZ'},
    @{Row=85; Col=4; Value=42},
    @{Row=86; Col=1; Value='08:01:10'},
    @{Row=86; Col=2; Value='08:33'},
    @{Row=86; Col=3; Value='10_OLMOS'},
    @{Row=86; Col=4; Value=32},
    @{Row=87; Col=1; Value='08:01:10'},
    @{Row=87; Col=2; Value='08:35'},
    @{Row=87; Col=3; Value='23_HERNANDEZ'},
    @{Row=87; Col=4; Value=34},
    @{Row=88; Col=2; Value='08:37'},
    @{Row=88; Col=3; Value='26_HERNANDEZ'},
    @{Row=88; Col=4; Value=49},
    @{Row=89; Col=1; Value='08:01:10'},
    @{Row=89; Col=2; Value='08:40'},
    @{Row=89; Col=3; Value='16_SANTA ANA'},
    @{Row=89; Col=4; Value=39},
    @{Row=90; Col=1; Value='07:17:59'},
    @{Row=90; Col=2; Value='08:41'},
    @{Row=90; Col=3; Value='81_EL PELIGRO'},
    @{Row=90; Col=4; Value=84},
    @{Row=91; Col=1; Value='08:01:10'},
    @{Row=91; Col=2; Value='08:42'},
    @{Row=91; Col=3; Value='81_EL PELIGRO'},
    @{Row=91; Col=4; Value=41},
    @{Row=92; Col=1; Value='07:17:59'},
    @{Row=92; Col=2; Value='08:43'},
    @{Row=92; Col=3; Value='14_ABASTO'},
    @{Row=92; Col=4; Value=86},
    @{Row=93; Col=1; Value='08:01:10'},
    @{Row=93; Col=2; Value='08:44'},
    @{Row=93; Col=3; Value='14_ABASTO'},
    @{Row=93; Col=4; Value=43},
    @{Row=94; Col=1; Value='08:01:10'},
    @{Row=94; Col=2; Value='08:49'},
    @{Row=94; Col=3; Value='26_HERNANDEZ'},
    @{Row=94; Col=4; Value=48},
    @{Row=95; Col=2; Value='08:53'},
    @{Row=95; Col=3; Value='17_ROMERO'},
    @{Row=95; Col=4; Value=96},
    @{Row=96; Col=1; Value='08:01:10'},
    @{Row=96; Col=2; Value='08:54'},
    @{Row=96; Col=3; Value='17_ROMERO'},
    @{Row=96; Col=4; Value=53},
    @{Row=97; Col=2; Value='09:01'},
    @{Row=97; Col=3; Value='215A_EL PATO'},
    @{Row=97; Col=4; Value=104},
    @{Row=98; Col=1; Value='08:01:10'},
    @{Row=98; Col=2; Value='09:02'},
    @{Row=98; Col=3; Value='215A_EL PATO'},
    @{Row=98; Col=4; Value=61},
    @{Row=99; Col=1; Value='08:01:10'},
    @{Row=99; Col=2; Value='09:04'},
    @{Row=99; Col=3; Value='11_ETCHEVERRY'},
    @{Row=99; Col=4; Value=63},
    @{Row=100; Col=1; Value='07:17:59'},
    @{Row=100; Col=2; Value='09:10'},
    @{Row=100; Col=3; Value='16_P MOR-SANTA ANA'},
    @{Row=100; Col=4; Value=113},
    @{Row=101; Col=1; Value='08:01:10'},
    @{Row=101; Col=2; Value='09:11'},
    @{Row=101; Col=3; Value='16_P MOR-SANTA ANA'},
    @{Row=101; Col=4; Value=70},
    @{Row=102; Col=1; Value='07:17:59'},
    @{Row=102; Col=2; Value='09:16'},
    @{Row=102; Col=3; Value='27_EL RETIRO'},
    @{Row=102; Col=4; Value=119},
    @{Row=103; Col=1; Value='08:01:10'},
    @{Row=103; Col=2; Value='09:17'},
    @{Row=103; Col=3; Value='27_EL RETIRO'},
    @{Row=103; Col=4; Value=76},
    @{Row=104; Col=1; Value='08:01:10'},
    @{Row=104; Col=2; Value='09:21'},
    @{Row=104; Col=3; Value='26_HERNANDEZ'},
    @{Row=104; Col=4; Value=80},
    @{Row=105; Col=1; Value='08:01:10'},
    @{Row=105; Col=2; Value='09:23'},
    @{Row=105; Col=3; Value='17_ROMERO'},
    @{Row=105; Col=4; Value=82},
    @{Row=106; Col=1; Value='08:01:10'},
    @{Row=106; Col=2; Value='09:24'},
    @{Row=106; Col=3; Value='11_ETCHEVERRY'},
    @{Row=106; Col=4; Value=83},
    @{Row=106; Col=5; Value='LP1912'},
    @{Row=107; Col=1; Value='07:48:05'},
    @{Row=107; Col=2; Value='09:32'},
    @{Row=107; Col=3; Value='15_ABASTO'},
    @{Row=107; Col=4; Value=104},
    @{Row=107; Col=5; Value='LP1912'},
    @{Row=108; Col=1; Value='08:01:10'},
    @{Row=108; Col=2; Value='09:33'},
    @{Row=108; Col=3; Value='10_OLMOS'},
    @{Row=108; Col=4; Value=92},
    @{Row=108; Col=5; Value='LP1912'},
    @{Row=109; Col=1; Value='07:48:05'},
    @{Row=109; Col=2; Value='09:34'},
    @{Row=109; Col=3; Value='16_SANTA ANA'},
    @{Row=109; Col=4; Value=106},
    @{Row=109; Col=5; Value='LP1912'},
    @{Row=110; Col=1; Value='08:01:10'},
    @{Row=110; Col=2; Value='09:42'},
    @{Row=110; Col=3; Value='215C_EL PATO'},
    @{Row=110; Col=4; Value=101},
    @{Row=110; Col=5; Value='LP1912'},
    @{Row=111; Col=1; Value='08:01:10'},
    @{Row=111; Col=2; Value='09:52'},
    @{Row=111; Col=3; Value='15_ABASTO'},
    @{Row=111; Col=4; Value=111},
    @{Row=111; Col=5; Value='LP1912'}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item('LP1912-215')
$updates = @(
    @{Row=2; Col=1; Value='Última actualización: 08:01:10'},
    @{Row=20; Col=1; Value='08:01:10'},
    @{Row=20; Col=4; Value=22},
    @{Row=22; Col=1; Value='08:01:10'},
    @{Row=22; Col=4; Value=61},
    @{Row=23; Col=1; Value='08:01:10'},
    @{Row=23; Col=4; Value=101}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item('6203-6173')
$updates = @(
    @{Row=2; Col=1; Value='Última actualización: 08:01:10'},
    @{Row=3; Col=1; Value='Total filas: 18'},
    @{Row=18; Col=1; Value='08:01:10'},
    @{Row=18; Col=2; Value='08:22'},
    @{Row=18; Col=4; Value=21},
    @{Row=19; Col=1; Value='07:48:05'},
    @{Row=19; Col=2; Value='08:25'},
    @{Row=19; Col=3; Value='215C_LA PLATA'},
    @{Row=19; Col=4; Value=37},
    @{Row=19; Col=5; Value='L6203'},
    @{Row=20; Col=1; Value='07:17:59'},
    @{Row=20; Col=2; Value='08:34'},
    @{Row=20; Col=4; Value=77},
    @{Row=21; Col=1; Value='08:01:10'},
    @{Row=21; Col=2; Value='08:35'},
    @{Row=21; Col=3; Value='215A_LA PLATA'},
    @{Row=21; Col=4; Value=34},
    @{Row=21; Col=5; Value='L6173'},
    @{Row=22; Col=1; Value='07:17:59'},
    @{Row=22; Col=2; Value='09:08'},
    @{Row=22; Col=4; Value=111},
    @{Row=23; Col=1; Value='08:01:10'},
    @{Row=23; Col=2; Value='09:09'},
    @{Row=23; Col=3; Value='215D_LA PLATA'},
    @{Row=23; Col=4; Value=68},
    @{Row=23; Col=5; Value='L6203'}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
